$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- "Menu" table (columns G/H) gets a new "price" (double) column,
#     and a new "foodExchangeOpen" (datetime) row is inserted before the
#     existing "mensaOpen" (bit) row, which moves down from row 7 to row 9. ---

# Row 7: was mensaOpen/bit -> becomes price/double
$ws.Range("G7").Value = "price"
$ws.Range("H7").Value = "double"

# Row 8 (new): foodExchangeOpen/datetime
$ws.Range("G8").Value = "foodExchangeOpen"
$ws.Range("H8").Value = "datetime"

# Row 9 (new): re-add the original mensaOpen/bit entry
$ws.Range("G9").Value = "mensaOpen"
$ws.Range("H9").Value = "bit"

# Copy the label formatting (bold-less, centered) from G7 onto the two
# newly added label cells G8:G9 so they match the rest of the column.
$ws.Range("G7").Copy()
$ws.Range("G8:G9").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- Column width / autofit tweaks ---
$ws.Columns("B:B").ColumnWidth = 10.166666666666666
$ws.Columns("G:G").ColumnWidth = 17.666666666666668
$ws.Columns("L:L").ColumnWidth = 15.5

# --- Selection moved to L8 ---
$ws.Range("L8").Select()

Write-Output "done"
